{"js": "// The document re-saves after a spelling/grammar pass and a cursor\n// (re-)position: some paragraphs get their runs split around new\n// <w:proofErr> markers (spell/grammar squiggles) that carry no text\n// of their own, the Word \"_GoBack\" bookmark moves from the\n// \"FULL ALIGNMENT...\" paragraph to the \"Have better buttons created\"\n// paragraph, and four blank paragraphs are appended at the very end\n// of the body (after \"CLEAN UP CODE!!\"). None of this changes the\n// visible text of the document.\n\nfunction flatOpcBody(innerXml) {\n  // Minimal single-part Flat OPC wrapper accepted by Range.insertOoxml /\n  // Body.insertOoxml - only the document.xml part is needed for a body\n  // fragment insertion.\n  return (\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + innerXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\n// Replace the single paragraph whose current text equals `oldText` with the\n// literal paragraph XML in `newParaXml` (a single <w:p>...</w:p>).\nasync function replaceParagraphByText(context, paragraphs, oldText, newParaXml) {\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === oldText) {\n      paragraphs.items[i].getRange().insertOoxml(flatOpcBody(newParaXml), Word.InsertLocation.replace);\n      return true;\n    }\n  }\n  throw new Error(\"Paragraph not found: \" + oldText);\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) \"Properly check into GitHub\" -> split off \"GitHub\" with a spell-check\n//    proofing-error wrapper (text unchanged).\nawait replaceParagraphByText(\n  context,\n  paragraphs,\n  \"Properly check into GitHub\",\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Properly check into </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>GitHub</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n  '</w:p>'\n);\n\n// 2) \"Have better buttons created\" -> the \"_GoBack\" bookmark (last edit\n//    position) now sits inside this paragraph, splitting the run.\nawait replaceParagraphByText(\n  context,\n  paragraphs,\n  \"Have better buttons created\",\n  '<w:p>' +\n    '<w:r><w:t>Have be</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r><w:t>tter buttons created</w:t></w:r>' +\n  '</w:p>'\n);\n\n// 3) \"Fix text color and alignment of  Navigation Controller text\" -> a\n//    grammar-check wrapper around \"of  Navigation\" (text unchanged).\nawait replaceParagraphByText(\n  context,\n  paragraphs,\n  \"Fix text color and alignment of  Navigation Controller text\",\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Fix text color and alignment </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>of  Navigation</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Controller text</w:t></w:r>' +\n  '</w:p>'\n);\n\n// 4) \"Have better VECTOR social buttons created\" -> the trailing run is\n//    split so a grammar-check wrapper surrounds just \"created\".\nawait replaceParagraphByText(\n  context,\n  paragraphs,\n  \"Have better VECTOR social buttons created\",\n  '<w:p>' +\n    '<w:r><w:t>Have better</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> VECTOR</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> social buttons </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>created</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n  '</w:p>'\n);\n\n// 5) \"FULL ALIGNMENT ON PAYMENT SCREEN\" -> loses the \"_GoBack\" bookmark\n//    that moved to paragraph (2) above; text unchanged.\nawait replaceParagraphByText(\n  context,\n  paragraphs,\n  \"FULL ALIGNMENT ON PAYMENT SCREEN\",\n  '<w:p><w:r><w:t>FULL ALIGNMENT ON PAYMENT SCREEN</w:t></w:r></w:p>'\n);\n\n// 6) Four new blank paragraphs appended at the very end of the body,\n//    right after \"CLEAN UP CODE!!\".\nlet lastIdx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"CLEAN UP CODE!!\") { lastIdx = i; }\n}\nif (lastIdx === -1) throw new Error('Paragraph not found: \"CLEAN UP CODE!!\"');\nparagraphs.items[lastIdx].getRange().insertOoxml(\n  flatOpcBody('<w:p/><w:p/><w:p/><w:p/>'),\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# The document re-saves after a spelling/grammar pass and a cursor\n# (re-)position: some paragraphs get their runs split around new\n# <w:proofErr> markers (spell/grammar squiggles) that carry no text\n# of their own, the Word \"_GoBack\" bookmark moves from the\n# \"FULL ALIGNMENT...\" paragraph to the \"Have better buttons created\"\n# paragraph, and four blank paragraphs are appended at the very end\n# of the body (after \"CLEAN UP CODE!!\"). None of this changes the\n# visible text of the document.\n\nfunction Get-FlatOpcBody {\n    param([string]$InnerXml)\n    return '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $InnerXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\nfunction Replace-ParagraphByText {\n    param($Doc, [string]$OldText, [string]$NewParaXml)\n    foreach ($p in $Doc.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $OldText) {\n            $p.Range.InsertXML((Get-FlatOpcBody $NewParaXml))\n            return\n        }\n    }\n    throw \"Paragraph not found: $OldText\"\n}\n\n$d = $word.ActiveDocument\n\n# 1) \"Properly check into GitHub\" -> split off \"GitHub\" with a spell-check\n#    proofing-error wrapper (text unchanged).\nReplace-ParagraphByText $d \"Properly check into GitHub\" (\n    '<w:p>' +\n        '<w:r><w:t xml:space=\"preserve\">Properly check into </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>GitHub</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p>'\n)\n\n# 2) \"Have better buttons created\" -> the \"_GoBack\" bookmark (last edit\n#    position) now sits inside this paragraph, splitting the run.\nReplace-ParagraphByText $d \"Have better buttons created\" (\n    '<w:p>' +\n        '<w:r><w:t>Have be</w:t></w:r>' +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n        '<w:bookmarkEnd w:id=\"0\"/>' +\n        '<w:r><w:t>tter buttons created</w:t></w:r>' +\n    '</w:p>'\n)\n\n# 3) \"Fix text color and alignment of  Navigation Controller text\" -> a\n#    grammar-check wrapper around \"of  Navigation\" (text unchanged).\nReplace-ParagraphByText $d \"Fix text color and alignment of  Navigation Controller text\" (\n    '<w:p>' +\n        '<w:r><w:t xml:space=\"preserve\">Fix text color and alignment </w:t></w:r>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:t>of  Navigation</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> Controller text</w:t></w:r>' +\n    '</w:p>'\n)\n\n# 4) \"Have better VECTOR social buttons created\" -> the trailing run is\n#    split so a grammar-check wrapper surrounds just \"created\".\nReplace-ParagraphByText $d \"Have better VECTOR social buttons created\" (\n    '<w:p>' +\n        '<w:r><w:t>Have better</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> VECTOR</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> social buttons </w:t></w:r>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:t>created</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n    '</w:p>'\n)\n\n# 5) \"FULL ALIGNMENT ON PAYMENT SCREEN\" -> loses the \"_GoBack\" bookmark\n#    that moved to paragraph (2) above; text unchanged.\nReplace-ParagraphByText $d \"FULL ALIGNMENT ON PAYMENT SCREEN\" (\n    '<w:p><w:r><w:t>FULL ALIGNMENT ON PAYMENT SCREEN</w:t></w:r></w:p>'\n)\n\n# 6) Four new blank paragraphs appended at the very end of the body,\n#    right after \"CLEAN UP CODE!!\".\n$endRange = $d.Content\n$endRange.Collapse(0)  # wdCollapseEnd\n$endRange.InsertXML((Get-FlatOpcBody '<w:p/><w:p/><w:p/><w:p/>'))\n"}
